# Remove the trailing "blank / page-break / copyright" paragraphs that
# followed the last "Requisitos" entry (LOM3013: Ciência dos Materiais),
# while leaving the final blank + page-break paragraphs (and the
# following section) untouched.

$d = $word.ActiveDocument

$reqIndex = -1
$copyIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*LOM3013: Ci*ncia dos Materiais (Requisito)*") {
        $reqIndex = $i
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $copyIndex = $i
    }
}

if ($reqIndex -lt 0) {
    throw "Could not locate the 'LOM3013: Ciencia dos Materiais (Requisito)' paragraph"
}
if ($copyIndex -lt 0) {
    throw "Could not locate the copyright paragraph"
}
if ($copyIndex -le $reqIndex) {
    throw "Unexpected paragraph ordering (req=$reqIndex copy=$copyIndex)"
}

# Delete every paragraph between (exclusive) the requirement line and
# (inclusive) the copyright paragraph: this removes the intervening
# blank paragraph, the page-break paragraph, and the copyright paragraph
# itself, in one shot.
$deleteStart = $d.Paragraphs.Item($reqIndex + 1).Range.Start
$deleteEnd = $d.Paragraphs.Item($copyIndex).Range.End

$r = $d.Range($deleteStart, $deleteEnd)
$r.Delete()

Write-Output "Removed paragraphs $($reqIndex + 1) through $copyIndex (inclusive)."
